# Adjusted Qty On Hand column after taking inventory.
# Digi-Key order has been placed.  They haven't given me a final total yet,
# so I still cannot compute the final kit price, but the quantity price
# breaks probably saved about $6 per kit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "On Hand" (column N) quantities discovered during inventory.
# All other formula columns (O/P/Q = Qty Needed / Purchase Units / Purchase
# Cost) as well as the section and grand-total rows recompute automatically.
$ws.Range("N29").Value = 10   # 2.1mmx5.5mm DC barrel jack
$ws.Range("N30").Value = 10   # audio jack
$ws.Range("N37").Value = 10   # push button switch
$ws.Range("N41").Value = 10   # shunt jumper
$ws.Range("N43").Value = 10   # 0.01uF capacitor

# Slightly widen column Q so the recalculated totals still fit comfortably.
$ws.Range("Q1").ColumnWidth = 12

# Restore the active cell left over from the editing session.
$ws.Range("G27").Select() | Out-Null
